$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "1790"
$t.Cell(5,1).Range.Text = "0.00001"
$t.Cell(6,1).Range.Text = "0.00258"
$t.Cell(7,1).Range.Text = "0.00014"
$t.Cell(8,1).Range.Text = "0.00009"
$t.Cell(9,1).Range.Text = "0.00021"
$t.Cell(10,1).Range.Text = "0.00024"
$t.Cell(11,1).Range.Text = "0.00031"
$t.Cell(12,1).Range.Text = "0.28270"

$t.Cell(44,1).Range.Text = "99.96"
$t.Cell(45,1).Range.Text = "0.28"
$t.Cell(46,1).Range.Text = "795"
